# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# on the active worksheet to reflect the latest scraped values.
#
# Column D ("Price") values that look like plain decimal numbers are
# written with a leading apostrophe so Excel stores them as text (matching
# the source data, which intentionally keeps prices such as "1.016" or
# "19.57" as text strings rather than numeric values), then the cell style
# is reset back to "Normal" so no extra number-format/quote-prefix style is
# left behind on the cell.
#
# Column E ("Volume(1h)") values already contain non-numeric characters
# (spaces and a percent sign) so Excel stores them as text automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.162.99"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").Value = "'1.651.17"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("D4").Value = "'1.016"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.88%  '
$ws.Range("D5").Value = "'216.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.63%  '
$ws.Range("D6").Value = "'0.5062"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.91%  '
$ws.Range("D7").Value = "'1.015"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").Value = "'0.2593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("D9").Value = "'0.06467"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.68%  '
$ws.Range("D10").Value = "'19.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("D11").Value = "'0.07776"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.27%  '
$ws.Range("D12").Value = "'1.668.05"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.74%  '
$ws.Range("D13").Value = "'4.280"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.53%  '
$ws.Range("D14").Value = "'1.876.35"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.48%  '
$ws.Range("D15").Value = "'0.5498"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.77%  '
$ws.Range("D16").Value = "'0.0₅7968"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.70%  '
$ws.Range("D17").Value = "'63.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.63%  '
$ws.Range("D18").Value = "'26.165.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.16%  '
$ws.Range("D19").Value = "'1.015"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.37%  '
$ws.Range("D20").Value = "'204.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.64%  '
$ws.Range("D21").Value = "'4.339"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.41%  '
$ws.Range("E22").Value = '  -0.40%  '
$ws.Range("D23").Value = "'5.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.14%  '
$ws.Range("D24").Value = "'1.015"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.56%  '
$ws.Range("D25").Value = "'1.963"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +11.83%  '
$ws.Range("D26").Value = "'142.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E27").Value = '  +0.33%  '
$ws.Range("D28").Value = "'15.87"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.54%  '
$ws.Range("D29").Value = "'6.780"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.59%  '
$ws.Range("D30").Value = "'0.05092"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("D32").Value = "'3.282"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.46%  '
$ws.Range("D33").Value = "'3.218"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.06%  '
$ws.Range("D34").Value = "'1.556"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = "'0.9027"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.31%  '
$ws.Range("D37").Value = "'2.639"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.98%  '
$ws.Range("E38").Value = '  -0.67%  '
$ws.Range("D39").Value = "'1.158.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.36%  '
$ws.Range("D40").Value = "'0.01582"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").Value = "'2.583"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("E42").Value = '  +0.51%  '
$ws.Range("D43").Value = "'5.697"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.69%  '
$ws.Range("D44").Value = "'0.8213"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.64%  '
$ws.Range("D45").Value = "'100.21"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("D46").Value = "'1.782.72"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.68%  '
$ws.Range("E47").Value = '  +4.33%  '
$ws.Range("D48").Value = "'0.4564"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.46%  '
$ws.Range("D49").Value = "'1.015"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("D50").Value = "'55.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.15%  '
$ws.Range("D51").Value = "'0.05057"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.71%  '
